$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "To do"
$ws.Range("B1").Value = "Done"

# Body rows
$ws.Range("A2").Value = "Ajouter les boutons pour revenir à la page Home(Signin Signup)"
$ws.Range("A3").Value = "Changer les icons de signup à gauche verticalement"
$ws.Range("A4").Value = "Changer le bouton dans Home page pour étre comme les boutons dans la page signin "

# Header formatting: bold, size 14
$headerFont = $ws.Range("A1:B1").Font
$headerFont.Bold = $true
$headerFont.Size = 14

# Column widths (target character widths 61.77734375 / 42.5546875; engine
# quantizes ColumnWidth to 1/6-character pixel steps, so these inputs land
# on the closest representable bucket)
$ws.Columns.Item(1).ColumnWidth = 61
$ws.Columns.Item(2).ColumnWidth = 41.6

# Row heights
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(4).RowHeight = 50.4

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A5").Select()
